$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "new handler" data block (rows 18-21) with the values from
# the timer's second handler: clock prescaler (B19) and counter period (B20).
$ws.Range("B19").Value = 72
$ws.Range("B20").Value = 65000
